$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B so its stored width matches column A (target 15.42578125).
# The ColumnWidth setter snaps to a pixel grid (steps of 1/6 character), so
# 14.7 is the input that lands on the closest achievable stored width (15.5).
$ws.Columns.Item(2).ColumnWidth = 14.7

# Update the numeric values in columns A and B (rows 1-32)
$ws.Range("A1").Value = -0.082199919036966662
$ws.Range("B1").Value = 0.082044355511740719
$ws.Range("A2").Value = -0.017998870150336188
$ws.Range("B2").Value = 0.017359958747904969
$ws.Range("A3").Value = 0.085573427404558089
$ws.Range("B3").Value = -0.085939011713129787
$ws.Range("A4").Value = -0.19805441659352141
$ws.Range("B4").Value = 0.1969462285217034
$ws.Range("A5").Value = -0.19094622889020929
$ws.Range("B5").Value = 0.18870814135191338
$ws.Range("A6").Value = -0.090163353240932942
$ws.Range("B6").Value = 0.090062455497230687
$ws.Range("A7").Value = -0.070062455950145264
$ws.Range("B7").Value = 0.069834691966160989
$ws.Range("A8").Value = -0.049834692423151417
$ws.Range("B8").Value = 0.049661399492535452
$ws.Range("A9").Value = -0.04366139988536144
$ws.Range("B9").Value = 0.043521783234600875
$ws.Range("A10").Value = -0.037521783631440542
$ws.Range("B10").Value = 0.037502921013988555
$ws.Range("A11").Value = -0.033002921403546992
$ws.Range("B11").Value = 0.032971524370594807
$ws.Range("A12").Value = -0.027391179744002514
$ws.Range("B12").Value = 0.027284336736391523
$ws.Range("A13").Value = -0.013795362012025691
$ws.Range("B13").Value = 0.013787851045753641
$ws.Range("A14").Value = -0.0017878514772613485
$ws.Range("B14").Value = 0.0017876422239329059
$ws.Range("A15").Value = 0.0042123573739312192
$ws.Range("B15").Value = -0.0042137049342230171
$ws.Range("A16").Value = -0.015027111459101405
$ws.Range("B16").Value = 0.015004368130623558
$ws.Range("A17").Value = -0.0090043685342235946
$ws.Range("B17").Value = 0.0089999995812855005
$ws.Range("A18").Value = -0.054045485744218524
$ws.Range("B18").Value = 0.05401826132517229
$ws.Range("A19").Value = -0.045018261703748674
$ws.Range("B19").Value = 0.044824239088566031
$ws.Range("A20").Value = -0.035824239472882269
$ws.Range("B20").Value = 0.035789143154836722
$ws.Range("A21").Value = -0.0090042868744601989
$ws.Range("B21").Value = 0.0089999996142031691
$ws.Range("A22").Value = -0.093933386531432461
$ws.Range("B22").Value = 0.093624992878691415
$ws.Range("A23").Value = -0.084624993261002146
$ws.Range("B23").Value = 0.084124937069428363
$ws.Range("A24").Value = -0.042124937619985481
$ws.Range("B24").Value = 0.041999999446503011
$ws.Range("A25").Value = -0.092544788493686525
$ws.Range("B25").Value = 0.092378103348984553
$ws.Range("A26").Value = -0.086378103733029121
$ws.Range("B26").Value = 0.086168100559024907
$ws.Range("A27").Value = -0.080168100944937315
$ws.Range("B27").Value = 0.079467597552962843
$ws.Range("A28").Value = -0.073467597946550889
$ws.Range("B28").Value = 0.072997930377275644
$ws.Range("A29").Value = -0.060997930805784861
$ws.Range("B29").Value = 0.060783822944078025
$ws.Range("A30").Value = -0.04078382341439335
$ws.Range("B30").Value = 0.040666861429760726
$ws.Range("A31").Value = -0.027019571732951775
$ws.Range("B31").Value = 0.027000839068552551
$ws.Range("A32").Value = -0.0060008395481849774
$ws.Range("B32").Value = 0.0059999995938486705
